$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 873 (shifts existing rows 873-914 down to 874-915)
$ws.Rows.Item(873).Insert()

# Populate the newly inserted row with the new entry.
# Force column A to be treated as text so the date-like string isn't
# auto-converted to a date serial number, then clear the formatting so
# no extra style index gets attached to the cell (matches surrounding rows).
$ws.Range("A873").NumberFormat = "@"
$ws.Cells.Item(873, 1).Value = "2026/02/28"
$ws.Cells.Item(873, 2).Value = "土"
$ws.Cells.Item(873, 3).Value = 1
$ws.Cells.Item(873, 4).Value = 24
$ws.Range("A873").ClearFormats()
